$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds monthly data for 2014-2017, 12 rows per year, starting at row 2
# (row 1 is the header). Each year's 12-row block must be reordered so that the
# October/November/December rows move to the top of the block, followed by the
# original January-September rows (i.e. a left-rotation by 9 / right-rotation by 3).
$numCols = 9
$yearStartRows = @(2, 14, 26, 38)
$rowsPerYear = 12

# 1) Snapshot every data row (rows 2-49, columns A-I) before we overwrite anything,
#    so reads are never affected by writes happening earlier in the loop.
$snapshot = @{}
for ($r = 2; $r -le 49; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) For each year block, figure out which old row feeds each new row position:
#    new offsets 0,1,2 <- old offsets 9,10,11 (Oct,Nov,Dec)
#    new offsets 3..11 <- old offsets 0..8   (Jan..Sep)
foreach ($start in $yearStartRows) {
    $oldOrder = @($start + 9, $start + 10, $start + 11)
    for ($i = 0; $i -lt 9; $i++) {
        $oldOrder += ,($start + $i)
    }

    for ($offset = 0; $offset -lt $rowsPerYear; $offset++) {
        $newRow = $start + $offset
        $oldRow = $oldOrder[$offset]
        $vals = $snapshot[$oldRow]
        for ($c = 1; $c -le $numCols; $c++) {
            $newVal = $vals[$c - 1]
            $cell = $ws.Cells.Item($newRow, $c)
            $curVal = $cell.Value()
            # Only write when the value actually changes, so cells that
            # legitimately stay blank/unchanged are left untouched.
            if ($curVal -ne $newVal) {
                $cell.Value = $newVal
            }
        }
    }
}
